$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - name unchanged, values updated
$ws.Range("B3").Value = 0.9985720097232855
$ws.Range("C3").Value = 0.9986063153366346
$ws.Range("D3").Value = 0.9864793705643443

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9981685101656613
$ws.Range("C4").Value = 0.9983793906161081
$ws.Range("D4").Value = 0.9689077937458949

# Row 5: AdaBoostRegressor -> MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9991267954765871
$ws.Range("C5").Value = 0.9990471027046756
$ws.Range("D5").Value = 0.9981279713680499
